$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell comment currently anchored at AA2 needs to end up on the cell
# that will become Z2 once the STAT column is removed below. Capture its
# text now and re-create it at the final location after the column
# deletions (deleting columns does not relocate existing comments).
$oldComment = $ws.Range("AA2").Comment
$commentText = $oldComment.Text()
$oldComment.Delete()

# Remove the AC column (CHOICE_TIER) first so deleting it doesn't shift
# the position of column Z (STAT), which is deleted next.
$ws.Range("AC1").EntireColumn.Delete()
$ws.Range("Z1").EntireColumn.Delete()

$ws.Range("Z2").AddComment($commentText) | Out-Null
